$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A15: was added with the "date only" style (s=3); normalize it to the
# same datetime style (s=2) used by every other completed row.
$ws.Range("A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new transfer-control row (row 16).
$ws.Range("A16").Value = 45856
$ws.Range("A16").NumberFormat = "YYYY-MM-DD"

$ws.Range("B16").Value = "diegoarrumou"
$ws.Range("C16").Value = "diegoarrumou"
$ws.Range("D16").Value = "2025-07-18 14:35:32"
$ws.Range("E16").Value = "2025-07-18 14:35:32"
$ws.Range("F16").Value = "2025-07-18 14:35:34"
$ws.Range("G16").Value = "2025-07-18 14:35:34"
$ws.Range("H16").Value = "2025-07-18 14:35:34"
$ws.Range("I16").Value = "2025-07-18 14:35:36"
$ws.Range("J16").Value = "2025-07-18 14:35:36"
$ws.Range("K16").Value = "0:00:00"
$ws.Range("L16").Value = "0:00:00"
$ws.Range("M16").Value = "0:00:04"
$ws.Range("N16").Value = "2025-07-18 14:35:37"
$ws.Range("O16").Value = "2025-07-18 14:35:38"
$ws.Range("P16").Value = "2025-07-18 14:35:40"
$ws.Range("Q16").Value = "2025-07-18 14:35:40"
$ws.Range("R16").Value = "2025-07-18 14:35:41"
$ws.Range("S16").Value = "0:00:00"
$ws.Range("T16").Value = "0:00:01"
$ws.Range("U16").Value = "0:00:04"
$ws.Range("V16").Value = "0:00:01"
$ws.Range("W16").Value = ""
